# Adds 4 additional test-case rows (3,4,5,6) to the "Datos" sheet of the
# generar_clave data-driven test workbook, mirroring the format already
# used by row 2, and fills row 2's previously-empty "resultadoEsperado"
# (E2) / "mensajeEsperado" (K2) cells.
#
# Column layout (row 1 headers):
#   A idCaso  B orientacion  C codigoTransaccion  D codigoError
#   E resultadoEsperado  F numeroDocumento  G usuario  H clave
#   I segundaClave  J tipoDocumento  K mensajeEsperado

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# --- propagate row 2's cell formatting down to the new rows 3-6 --------
$ws.Range("A2:K2").Copy()
$ws.Range("A3:K6").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

$tipoDoc = "C" + [char]0x00E9 + "dula de ciudadan" + [char]0x00ED + "a"
$msgInvalido = "Usuario o clave inv" + [char]0x00E1 + "lida. Int" + [char]0x00E9 + "ntalo nuevamente"
$msgActivar = " Debes activar la clave en la Sucursal F" + [char]0x00ED + "sica. Si eres colombiano en el exterior comun" + [char]0x00ED + "cate con la Sucursal Telef" + [char]0x00F3 + "nica."
$msgBloqueada = "La clave que usas en el cajero est" + [char]0x00E1 + " bloqueada. Debes activarla en la Sucursal F" + [char]0x00ED + "sica. Para mayor informaci" + [char]0x00F3 + "n comun" + [char]0x00ED + "cate con la Sucursal Telef" + [char]0x00F3 + "nica."

# --- row 2: fill the previously blank resultadoEsperado/mensajeEsperado,
#     and swap in the new numeroDocumento/usuario for this case
$ws.Range("E2").Value = $msgInvalido
$ws.Range("F2").Value = "1001945795"
$ws.Range("G2").Value = "chipote95"
$ws.Range("K2").Value = $msgInvalido

# --- row 3: new case 2 / Acierto ---------------------------------------
$ws.Range("A3").Value = "2"
$ws.Range("B3").Value = "Acierto"
$ws.Range("F3").Value = "1245123880"
$ws.Range("G3").Value = "OSVPPRU96"
$ws.Range("H3").Value = 1234
$ws.Range("I3").Value = 1234
$ws.Range("J3").Value = $tipoDoc
$ws.Range("K3").Value = $msgActivar

# --- row 4: new case 3 / Alterno ---------------------------------------
$ws.Range("A4").Value = "3"
$ws.Range("B4").Value = "Alterno"
$ws.Range("F4").Value = "25130110"
$ws.Range("G4").Value = "USUCDTC1"
$ws.Range("H4").Value = 4321
$ws.Range("I4").Value = 1234
$ws.Range("J4").Value = $tipoDoc
$ws.Range("K4").Value = $msgActivar

# --- row 5: new case 4 / Alterno (no numeroDocumento/usuario/claves) ---
$ws.Range("A5").Value = "4"
$ws.Range("B5").Value = "Alterno"
$ws.Range("E5").ClearFormats()
$ws.Range("E5").Value = $msgBloqueada
$ws.Range("K5").Value = $msgBloqueada

# --- row 6: new case 5 / Alterno (idCaso/orientacion only) -------------
$ws.Range("A6").Value = "5"
$ws.Range("B6").Value = "Alterno"

# --- data validation list now spans the whole new range ----------------
$ws.Range("B2:B6").Validation.Delete()
$ws.Range("B2:B6").Validation.Add(3, 1, 1, "Listas!`$A`$2:`$A`$3", "0")
